# Tuckshop workbook fix: Special Instructions are now printed on every
# label, so the "Special instructions" column (R) no longer needs to hold
# the full free-text note - it just needs a short flag value, and the
# column is widened so the flag / wrapped text are easy to read, with the
# flag cell left selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the long allergy note in R2 (Special instructions) with a short
# "SPECIAL" marker, since the detail is now printed on every label.
$ws.Range("R2").Value = "SPECIAL"

# Widen column R (18) so the new marker / wrapped instructions are legible.
$ws.Columns.Item(18).ColumnWidth = 52.7109375

# Leave the Special instructions cell selected/active.
$ws.Range("R2").Select()
